$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue($cellRef, $value) {
    $range = $ws.Range($cellRef)
    $range.NumberFormat = "@"
    $range.Value = $value
    $range.Style = "Normal"
}

Set-TextValue 'D2' '28.018.76'
Set-TextValue 'E2' '  -0.99%  '
Set-TextValue 'D3' '1.761.18'
Set-TextValue 'E3' '  -1.57%  '
Set-TextValue 'D4' '1.002'
Set-TextValue 'E4' '  +0.04%  '
Set-TextValue 'D5' '335.33'
Set-TextValue 'E5' '  -1.30%  '
Set-TextValue 'D6' '0.9998'
Set-TextValue 'E6' '  -0.04%  '
Set-TextValue 'D7' '0.3917'
Set-TextValue 'E7' '  +2.14%  '
Set-TextValue 'D8' '0.3402'
Set-TextValue 'E8' '  -1.32%  '
Set-TextValue 'D9' '45.33'
Set-TextValue 'E9' '  -3.55%  '
Set-TextValue 'D10' '1.121'
Set-TextValue 'E10' '  -3.07%  '
Set-TextValue 'D11' '0.07248'
Set-TextValue 'E11' '  -2.34%  '
Set-TextValue 'D12' '1.000'
Set-TextValue 'E12' '  +0.18%  '
Set-TextValue 'D13' '22.30'
Set-TextValue 'E13' '  -4.31%  '
Set-TextValue 'D14' '6.160'
Set-TextValue 'E14' '  -4.87%  '
Set-TextValue 'B15' 'Chainlink'
Set-TextValue 'C15' 'https://coinranking.com/coin/VLqpJwogdhHNb+chainlink-link'
Set-TextValue 'D15' '7.119'
Set-TextValue 'E15' '  -3.84%  '
Set-TextValue 'B16' 'WrappedEther'
Set-TextValue 'C16' 'https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth'
Set-TextValue 'D16' '1.759.49'
Set-TextValue 'E16' '  -1.34%  '
Set-TextValue 'D17' '0.00001059'
Set-TextValue 'E17' '  -2.00%  '
Set-TextValue 'D18' '0.06630'
Set-TextValue 'E18' '  -0.93%  '
Set-TextValue 'D19' '80.56'
Set-TextValue 'E19' '  -2.36%  '
Set-TextValue 'D20' '0.9991'
Set-TextValue 'E20' '  -0.11%  '
Set-TextValue 'D21' '16.96'
Set-TextValue 'E21' '  -3.45%  '
Set-TextValue 'D22' '6.230'
Set-TextValue 'E22' '  -4.00%  '
Set-TextValue 'D23' '28.008.50'
Set-TextValue 'E23' '  -0.96%  '
Set-TextValue 'D24' '11.65'
Set-TextValue 'E24' '  -3.87%  '
Set-TextValue 'D25' '2.385'
Set-TextValue 'E25' '  +1.05%  '
Set-TextValue 'D26' '154.98'
Set-TextValue 'E26' '  +0.34%  '
Set-TextValue 'D27' '19.98'
Set-TextValue 'E27' '  -3.94%  '
Set-TextValue 'D28' '2.320'
Set-TextValue 'E28' '  -4.60%  '
Set-TextValue 'D29' '1.959.40'
Set-TextValue 'E29' '  -1.21%  '
Set-TextValue 'D30' '1.276'
Set-TextValue 'E30' '  -12.16%  '
Set-TextValue 'D31' '129.27'
Set-TextValue 'E31' '  -5.03%  '
Set-TextValue 'D32' '4.077'
Set-TextValue 'E32' '  +3.07%  '
Set-TextValue 'D33' '5.828'
Set-TextValue 'E33' '  -5.35%  '
Set-TextValue 'D34' '0.08723'
Set-TextValue 'E34' '  -2.41%  '
Set-TextValue 'D35' '12.10'
Set-TextValue 'E35' '  -5.70%  '
Set-TextValue 'B36' 'Hedera'
Set-TextValue 'C36' 'https://coinranking.com/coin/jad286TjB+hedera-hbar'
Set-TextValue 'D36' '0.06178'
Set-TextValue 'E36' '  -3.46%  '
Set-TextValue 'B37' 'VeChain'
Set-TextValue 'C37' 'https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet'
Set-TextValue 'D37' '0.02289'
Set-TextValue 'E37' '  -5.95%  '
Set-TextValue 'D38' '5.142'
Set-TextValue 'E38' '  -4.11%  '
Set-TextValue 'D39' '0.6509'
Set-TextValue 'E39' '  -5.56%  '
Set-TextValue 'D40' '0.2111'
Set-TextValue 'E40' '  -3.11%  '
Set-TextValue 'D41' '1.498'
Set-TextValue 'E41' '  -0.08%  '
Set-TextValue 'E42' '  -3.61%  '
Set-TextValue 'D43' '0.9994'
Set-TextValue 'E43' '  -0.02%  '
Set-TextValue 'D44' '7.872'
Set-TextValue 'E44' '  -5.48%  '
Set-TextValue 'D45' '13.74'
Set-TextValue 'E45' '  -3.15%  '
Set-TextValue 'E46' '  -1.31%  '
Set-TextValue 'D47' '0.6008'
Set-TextValue 'E47' '  -5.05%  '
Set-TextValue 'D48' '127.01'
Set-TextValue 'E48' '  -5.06%  '
Set-TextValue 'D49' '1.999'
Set-TextValue 'E49' '  -4.43%  '
Set-TextValue 'D50' '1.159'
Set-TextValue 'E50' '  -4.45%  '
Set-TextValue 'D51' '0.07002'
Set-TextValue 'E51' '  -6.75%  '
